# Updates cryptos list (price / 1h volume change columns) as scraped on
# Fri Sep 15 15:50:15 UTC 2023 by GitHub Actions. Also reorders rows 35/36
# (Maker now ranks above HuobiToken) and refreshes their price/volume.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values look like plain numbers (e.g. "212.67"), so force
# the cell to Text format first -- otherwise Excel's smart-parsing would
# silently convert them into numeric cells and drop the original string
# formatting used throughout this sheet (inline/shared text, not numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.408.41"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.624.01"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.67"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.95"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.35"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.77"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.434.50"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.28"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.68"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.215.44"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.795"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.501"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.761.81"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.87"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.73"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -0.91%  "
